$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Course descriptions added in column D. Order chosen to match original
# authoring sequence (languages typed first, then math, then CS) so that
# the shared-string table indices line up with the source workbook.

$ws.Range("D70").Value = "Students develop the ability to communicate about Personal and Family Life, School Life, Social Life, and Community Life using simple sentences containing basic language structure. This course counts towards the world languages course for the Advanced Studies Diploma."
$ws.Range("D70").Style = "Normal"

$ws.Range("D71").Value = "Students learn to function in real-life situations using more complex language structures and a wider range of vocabulary. Explore themes of Home Life, Student Life, Leisure Time, and Vacation and Travel. "
$ws.Range("D71").Style = "Normal"

$ws.Range("D58").Value = "Students will learn to use simple sentence structures and basic language structures to discuss about Personal and Family Life, School Life, Social Life, and Community Life. They will develop their listening, speaking, reading, and writing skills."
$ws.Range("D58").Style = "Normal"

$ws.Range("D59").Value = "Students continue to develop their skills in listening, speaking, reading, and writing, while learning how to function in real-life situations using more complex sentences and language structures. They will real material on familiar topics and produce short writing samples. "
$ws.Range("D59").Style = "Normal"

$ws.Range("D60").Value = "Students continue to develop their skills in listening, speaking, reading, and writing, while learning how to use more complex language structures on more abstract concepts. The themes of the class include Rights and Responsibilities, Future Plans and Choices, Teen Culture, Environment, and Humanities."
$ws.Range("D60").Style = "Normal"

$ws.Range("D64").Value = "Students learn basic language structure and pronunciation in order to read simple Latin passages. The relationship of English to Latin is emphasized in vocabulary, word derivation, and meanings of prefixes and suffixes. Students also learn about the geography, history, government, and culture of the Roman Empire. "
$ws.Range("D64").Style = "Normal"

$ws.Range("D65").Value = "Students learn more vocabulary, more complex language structures and syntax so that they are able to read more challenging passages in Latin. Students continue to study Roman life and Rome’s contribution to our civilization. "
$ws.Range("D65").Style = "Normal"

$ws.Range("D66").Value = "Students develop and refine their reading skills, learn additional vocabulary, and learn more complex language structures and syntax. Through translation and interpretation, students gain a greater understanding of the foundation of Western government and civilization."
$ws.Range("D66").Style = "Normal"

$ws.Range("D28").Value = "Topics include linear equations and inequalities, systems of linear equations, relations, functions, polynomials, and statistics. Emphasis is placed on making connections in algebra to geometry and statistics."
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").Value = "This class goes beyond the scope of Algebra I. Students are expected to master algebraic mechanics and understand the underlying theory, as well as apply the concepts to real-world situations. Emphasis is placed on algebraic connections to arithmetic, geometry, and statistics. "
$ws.Range("D29").Style = "Normal"

$ws.Range("D36").Value = "This course emphasizes two- and three-dimensional reasoning skills, coordinate and transformational geometry, and the use of geometric models to solve problems. "
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").Value = "Goes beyond the scope of Geometry. Heavily uses proofs to verify theorems. Students investigate non-Euclidean geometries and formal logic."
$ws.Range("D37").Style = "Normal"
$ws.Range("D37").VerticalAlignment = -4108

$ws.Range("D34").Value = "Topics include function, polynomials, rational expressions, complex numbers, exponential and logarithmic equations, arithmetic and geometric sequences and series, and data analysis. "
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").Value = "Students are expected to not only master algebraic mechanics but also to understand the underlying theory and to apply concepts to real-world situations in a meaningful way. Additional topics include matrices, infinite geometric sequences and series, permutations and combinations, and selected topics in discrete math. Emphasis is on modeling, logic, and interpretation of related graphs. "
$ws.Range("D35").Style = "Normal"

$ws.Range("D43").Value = "Students learn how to code in Java, developing their skills in defining, writing, and running programs on a computer. Students will work with both mathematical and non-mathematical problems. "
$ws.Range("D43").Style = "Normal"

$ws.Range("D61").Value = "Students will learn to use simple sentence structures and basic language structures to discuss about Personal and Family Life, School Life, Social Life, and Community Life. They will develop their listening, speaking, reading, and writing skills."
$ws.Range("D61").Style = "Normal"

$ws.Range("D62").Value = "Students continue to develop their skills in listening, speaking, reading, and writing, while learning how to function in real-life situations using more complex sentences and language structures. They will real material on familiar topics and produce short writing samples. "
$ws.Range("D62").Style = "Normal"

$ws.Range("D63").Value = "Students continue to develop their skills in listening, speaking, reading, and writing, while learning how to use more complex language structures on more abstract concepts. The themes of the class include Rights and Responsibilities, Future Plans and Choices, Teen Culture, Environment, and Humanities."
$ws.Range("D63").Style = "Normal"

$ws.Range("D67").Value = "Students will learn to use simple sentence structures and basic language structures to discuss about Personal and Family Life, School Life, Social Life, and Community Life. They will develop their listening, speaking, reading, and writing skills."
$ws.Range("D67").Style = "Normal"

$ws.Range("D68").Value = "Students continue to develop their skills in listening, speaking, reading, and writing, while learning how to function in real-life situations using more complex sentences and language structures. They will real material on familiar topics and produce short writing samples. "
$ws.Range("D68").Style = "Normal"

$ws.Range("D69").Value = "Students continue to develop their skills in listening, speaking, reading, and writing, while learning how to use more complex language structures on more abstract concepts. The themes of the class include Rights and Responsibilities, Future Plans and Choices, Teen Culture, Environment, and Humanities."
$ws.Range("D69").Style = "Normal"

$ws.Range("D72").Value = "Students will learn to use simple sentence structures and basic language structures to discuss about Personal and Family Life, School Life, Social Life, and Community Life. They will develop their listening, speaking, reading, and writing skills."
$ws.Range("D72").Style = "Normal"

$ws.Range("D73").Value = "Students continue to develop their skills in listening, speaking, reading, and writing, while learning how to function in real-life situations using more complex sentences and language structures. They will real material on familiar topics and produce short writing samples. "
$ws.Range("D73").Style = "Normal"

# --- sheet view: selection + scroll state ---
$ws.Range("D39").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
